$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format all target cells as Text so that numeric-looking strings
# (prices and percentages) are preserved verbatim instead of being
# auto-converted into numbers by Excel.
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","E48","D49","E49","E50","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "316.20"
$ws.Range("E2").Value = "1.43%"
$ws.Range("D3").Value = "37.87"
$ws.Range("E3").Value = "1.06%"
$ws.Range("D4").Value = "5.183"
$ws.Range("E4").Value = "1.50%"
$ws.Range("D5").Value = "0.07994"
$ws.Range("E5").Value = "1.41%"
$ws.Range("D6").Value = "4.488"
$ws.Range("E6").Value = "1.93%"
$ws.Range("D7").Value = "8.525"
$ws.Range("E7").Value = "3.28%"
$ws.Range("D8").Value = "1.919"
$ws.Range("E8").Value = "0.67%"
$ws.Range("E9").Value = "0.36%"
$ws.Range("D10").Value = "0.9438"
$ws.Range("E10").Value = "2.40%"
$ws.Range("D11").Value = "0.1312"
$ws.Range("E11").Value = "9.24%"
$ws.Range("D12").Value = "0.1936"
$ws.Range("E12").Value = "0.46%"
$ws.Range("D13").Value = "0.09088"
$ws.Range("E13").Value = "0.69%"
$ws.Range("D14").Value = "0.03396"
$ws.Range("E14").Value = "2.29%"
$ws.Range("D15").Value = "0.09521"
$ws.Range("E15").Value = "-1.08%"
$ws.Range("D16").Value = "0.001400"
$ws.Range("E16").Value = "1.70%"
$ws.Range("D17").Value = "0.005984"
$ws.Range("E17").Value = "-3.06%"
$ws.Range("D18").Value = "3.437"
$ws.Range("E18").Value = "-3.01%"
$ws.Range("D19").Value = "0.3515"
$ws.Range("E19").Value = "2.13%"
$ws.Range("D20").Value = "6.570"
$ws.Range("E20").Value = "25.60%"
$ws.Range("E21").Value = "1.32%"
$ws.Range("D22").Value = "0.2420"
$ws.Range("E22").Value = "-6.49%"
$ws.Range("D23").Value = "0.04364"
$ws.Range("E23").Value = "0.38%"
$ws.Range("E24").Value = "-1.28%"
$ws.Range("D25").Value = "0.004267"
$ws.Range("E25").Value = "-8.82%"
$ws.Range("D26").Value = "0.0001329"
$ws.Range("E26").Value = "-2.07%"
$ws.Range("D27").Value = "0.0003980"
$ws.Range("E27").Value = "-0.03%"
$ws.Range("D39").Value = "0.02393"
$ws.Range("E39").Value = "5.36%"
$ws.Range("D40").Value = "0.05153"
$ws.Range("E40").Value = "1.49%"
$ws.Range("D41").Value = "0.007648"
$ws.Range("E41").Value = "2.34%"
$ws.Range("E42").Value = "3.04%"
$ws.Range("D43").Value = "0.008524"
$ws.Range("E43").Value = "-5.48%"
$ws.Range("D44").Value = "0.002040"
$ws.Range("E44").Value = "4.81%"
$ws.Range("D45").Value = "0.008691"
$ws.Range("E45").Value = "-6.33%"
$ws.Range("D46").Value = "0.00006463"
$ws.Range("E46").Value = "-1.22%"
$ws.Range("E47").Value = "0.13%"
$ws.Range("E48").Value = "-14.88%"
$ws.Range("D49").Value = "0.001685"
$ws.Range("E49").Value = "68.88%"
$ws.Range("E50").Value = "0.13%"
$ws.Range("E51").Value = "0.13%"
